$d = $word.ActiveDocument
$ns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# 1. Remove the "Meta description" paragraph that currently sits right after
#    the title. Its text is re-purposed below as the new closing paragraph,
#    without the "Meta description:" label.
$metaPara = $d.Paragraphs.Item(2)
if ($metaPara.Range.Text -like "Meta description*") {
    $metaPara.Range.Delete()
}

# 2. Replace the final paragraph (currently the italic AI image-generation
#    prompt) with two new paragraphs: a bold title line followed by the
#    italic review blurb that used to live in the "Meta description" line.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
if ($lastPara.Range.Text -like "Create an eye-catching feature image*") {
    $xml = "<w:p $ns><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Crusader Free: Review of the Medieval-Themed Slot</w:t></w:r></w:p>" +
           "<w:p $ns><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Read our review of Crusader, a medieval-themed online slot game. Play for free and trigger free spins with the wild symbol.</w:t></w:r></w:p>"
    [void]$lastPara.Range.InsertXML($xml)
}
